$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("H4").Value = 1

# Row 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Rows 6-18: set H column to 1
for ($r = 6; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
